$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "MCT-3A-Máquinas Térmicas e de Fluxo"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "[André Guimarães-Desenho Técnico-1A, -]"

$ws.Range("D3").Value = "[-, -]"

$ws.Range("B7").Value = "-"
$ws.Range("C7").Value = "['MCT-2A-CAD', 'MCT-2A-CAD']"
$ws.Range("F7").Value = "-"

$ws.Range("B8").Value = "MCT-3A-Máquinas Térmicas e de Fluxo"
$ws.Range("D8").Value = "-"

$ws.Range("C10").Value = "[-, -]"

$ws.Range("B12").Value = "[-, -]"

$ws.Range("D14").Value = "[-, -]"

$ws.Range("B18").Value = "[-, -]"
$ws.Range("D18").Value = "[-, -]"
$ws.Range("E18").Value = "[-, -]"

$ws.Range("C20").Value = "[-, -]"
$ws.Range("D20").Value = "[-, -]"
